$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing row 243 (old rows 243-256 shift down to 245-258)
$ws.Rows.Item(243).Insert()
$ws.Rows.Item(243).Insert()

# New row 243: Doctor Davis / Primera
$ws.Cells.Item(243, 1).Value2 = 7
$ws.Cells.Item(243, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(243, 3).Value2 = "Ñuble"
$ws.Cells.Item(243, 4).Value2 = 44610
$ws.Cells.Item(243, 5).Value2 = 16
$ws.Cells.Item(243, 6).Value2 = "Fruta"
$ws.Cells.Item(243, 7).Value2 = 100103
$ws.Cells.Item(243, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(243, 9).Value2 = 100103004
$ws.Cells.Item(243, 10).Value2 = "Durazno"
$ws.Cells.Item(243, 11).Value2 = "Doctor Davis"
$ws.Cells.Item(243, 12).Value2 = "Primera"
$ws.Cells.Item(243, 13).Value2 = 120
$ws.Cells.Item(243, 14).Value2 = 12000
$ws.Cells.Item(243, 15).Value2 = 13000
$ws.Cells.Item(243, 16).Value2 = 12500
$ws.Cells.Item(243, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(243, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(243, 19).Value2 = 781
$ws.Cells.Item(243, 20).Value2 = 16

# New row 244: Doctor Davis / Segunda
$ws.Cells.Item(244, 1).Value2 = 7
$ws.Cells.Item(244, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(244, 3).Value2 = "Ñuble"
$ws.Cells.Item(244, 4).Value2 = 44610
$ws.Cells.Item(244, 5).Value2 = 16
$ws.Cells.Item(244, 6).Value2 = "Fruta"
$ws.Cells.Item(244, 7).Value2 = 100103
$ws.Cells.Item(244, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(244, 9).Value2 = 100103004
$ws.Cells.Item(244, 10).Value2 = "Durazno"
$ws.Cells.Item(244, 11).Value2 = "Doctor Davis"
$ws.Cells.Item(244, 12).Value2 = "Segunda"
$ws.Cells.Item(244, 13).Value2 = 200
$ws.Cells.Item(244, 14).Value2 = 10000
$ws.Cells.Item(244, 15).Value2 = 11000
$ws.Cells.Item(244, 16).Value2 = 10500
$ws.Cells.Item(244, 17).Value2 = "$/caja 16 kilos empedrada"
$ws.Cells.Item(244, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(244, 19).Value2 = 656
$ws.Cells.Item(244, 20).Value2 = 16
